# Insert a new price-record row before the current row 80 ("Terminal La
# Palmera de La Serena" / "Ajo", weekly update), shifting the existing
# rows 80..194 down to 81..195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 80 downward (keeps row 1 header + rows 2-79 untouched).
$ws.Rows(80).Insert()

# Populate the newly inserted row 80 with the new observation.
$ws.Cells.Item(80, 1).Value = 8
$ws.Cells.Item(80, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(80, 3).Value = "Coquimbo"
$ws.Cells.Item(80, 4).Value = 44557
$ws.Cells.Item(80, 5).Value = 4
$ws.Cells.Item(80, 6).Value = 100112003
$ws.Cells.Item(80, 7).Value = "Ajo"
$ws.Cells.Item(80, 8).Value = "Chino"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 520
$ws.Cells.Item(80, 11).Value = 19000
$ws.Cells.Item(80, 12).Value = 20000
$ws.Cells.Item(80, 13).Value = 19500
$ws.Cells.Item(80, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(80, 15).Value = "China"
$ws.Cells.Item(80, 16).Value = 1950
$ws.Cells.Item(80, 17).Value = 10
$ws.Cells.Item(80, 18).Value = "Hortaliza"
